{"js": "// The document (WordDocuments/Aptos/0483.docx) is rewritten from a \"Quantum Computing\"\n// themed essay to a \"Perception\" themed essay, the author's name/email are changed, and\n// the essay body grows substantially (new sentences + a new \"Body:\" section). Because so\n// much of the running text is rewritten/expanded sentence-by-sentence, the body paragraph\n// and the summary paragraph are replaced wholesale with their final OOXML (preserving the\n// exact run/formatting conventions already used in the document: Aptos font, black color,\n// sz 24 for the body, default size for the summary, <w:br/> separated sub-paragraphs, etc.).\n// The short, uniquely-identifiable title / author-name / author-email runs are updated with\n// simple search & replace so that their original run formatting is left completely untouched.\n\n// 1) Title\n{\n  const results = context.document.body.search(\"Quantum Computing: A Journey into the Quantum Realm\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"The Art of Perception: Deconstructing the Fabric of Reality\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Author name\n{\n  const results = context.document.body.search(\"Samuel Grey\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"Emily Carter\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Author email (three separate runs: \"samuel\" / \"grey@quantumcomputing\" / \"net\")\n{\n  const results = context.document.body.search(\"samuel\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"emily\", Word.InsertLocation.replace);\n  await context.sync();\n}\n{\n  const results = context.document.body.search(\"grey@quantumcomputing\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"carter@highschool\", Word.InsertLocation.replace);\n  await context.sync();\n}\n{\n  const results = context.document.body.search(\"net\", { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(\"edu\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) Replace the whole introductory paragraph (the one containing the <w:br/> separated\n//    mini-paragraphs) with its new \"Perception\" themed content.\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const introParagraph = paras.items[4];\n  introParagraph.getRange().insertOoxml(\"<?xml version=\\\"1.0\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>The world we perceive is a symphony of sensations, a tapestry woven from the threads of our senses</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Our eyes, ears, nose, tongue, and skin act as gateways, translating physical stimuli into a subjective reality</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> This intricate interplay between sensory inputs and our brains gives rise to our perception of the world, shaping our understanding of reality</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Yet, this perception is not a passive process but rather an active engagement, a dance between our senses and our minds</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/><w:t>We perceive not merely with our senses but also with our minds, our memories, and our expectations</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Our past experiences, cultural background, and beliefs influence how we interpret sensory information</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> What we see, hear, smell, taste, and touch is filtered through the lens of our individual experiences, creating a unique and personal reality for each of us</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> This interplay between our senses, our minds, and our past experiences creates a dynamic and ever-changing perception of the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/><w:t>Perception is a complex process that defies simple definitions</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> It involves not only the physical reception of stimuli but also the cognitive interpretation of those stimuli</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> It is a process that is both subjective and objective, personal and universal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> It is the foundation of our understanding of the world and the key to unlocking the secrets of our consciousness</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/><w:t>Body:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/><w:t>Our perception is shaped by our senses</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> The eyes, ears, nose, tongue, and skin are the primary gateways through which we receive information about the world around us</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> These sensory organs convert physical stimuli into electrical signals that are transmitted to the brain for interpretation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> The brain then integrates these signals, creating a coherent and meaningful representation of the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> This process is incredibly complex and involves a multitude of neural pathways and brain regions working in concert</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:lastRenderedPageBreak/><w:br/><w:t>Perception is influenced by our past experiences</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Our memory plays a crucial role in shaping our perception of the present</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Past experiences create expectations about what we expect to see, hear, smell, taste, and touch</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> These expectations can influence how we interpret sensory information, leading to biased or inaccurate perceptions</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> For example, a person who has had a negative experience with a particular food may be more likely to perceive that food as being unappealing, even if it is not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:br/><w:t>Perception is also influenced by our beliefs and expectations</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Our cultural background, beliefs, and expectations can shape how we perceive the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> For example, a person who believes in a particular religious doctrine may interpret ambiguous sensory information in a way that confirms their beliefs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Similarly, a person who expects to see a certain outcome may be more likely to perceive evidence that supports that outcome, even if the evidence is weak or nonexistent</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/><w:sz w:val=\\\"24\\\"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 5) Replace the Summary paragraph's text with its new \"Perception\" themed content.\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const summaryParagraph = paras.items[paras.items.length - 1];\n  summaryParagraph.getRange().insertOoxml(\"<?xml version=\\\"1.0\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t>Perception is a complex and multifaceted process that involves the interaction of our senses, our minds, and our past experiences</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> It is a process that is both subjective and objective, personal and universal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> Perception is the foundation of our understanding of the world and the key to unlocking the secrets of our consciousness</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> By understanding the nature of perception, we can gain a deeper appreciation for the intricate beauty of our subjective realities and the common threads that bind us together as a human species</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Aptos\\\" w:hAnsi=\\\"Aptos\\\"/><w:color w:val=\\\"000000\\\"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 6) A new, empty trailing paragraph is added at the very end of the document body.\n{\n  context.document.body.insertParagraph(\"\", Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# The document (WordDocuments/Aptos/0483.docx) is rewritten from a \"Quantum Computing\"\n# themed essay to a \"Perception\" themed essay, the author's name/email are changed, and\n# the essay body grows substantially (new sentences + a new \"Body:\" section). Because so\n# much of the running text is rewritten/expanded sentence-by-sentence, the body paragraph\n# and the summary paragraph are replaced wholesale with their final WordOpenXML (preserving\n# the exact run/formatting conventions already used in the document: Aptos font, black\n# color, sz 24 for the body, default size for the summary, <w:br/> separated sub-paragraphs,\n# etc.). The short, uniquely identifiable title / author-name / author-email runs are\n# updated with simple Find & Replace so their original run formatting is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Unique($doc, $findText, $replaceText, $wholeWord) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Title\nReplace-Unique $d \"Quantum Computing: A Journey into the Quantum Realm\" \"The Art of Perception: Deconstructing the Fabric of Reality\" $false\n\n# 2) Author name\nReplace-Unique $d \"Samuel Grey\" \"Emily Carter\" $false\n\n# 3) Author email (three separate runs: \"samuel\" / \"grey@quantumcomputing\" / \"net\")\nReplace-Unique $d \"samuel\" \"emily\" $false\nReplace-Unique $d \"grey@quantumcomputing\" \"carter@highschool\" $false\nReplace-Unique $d \"net\" \"edu\" $true\n\n# 4) Replace the whole introductory paragraph (the one containing the <w:br/> separated\n#    mini-paragraphs) with its new \"Perception\" themed content.\n$introXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>The world we perceive is a symphony of sensations, a tapestry woven from the threads of our senses</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Our eyes, ears, nose, tongue, and skin act as gateways, translating physical stimuli into a subjective reality</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> This intricate interplay between sensory inputs and our brains gives rise to our perception of the world, shaping our understanding of reality</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Yet, this perception is not a passive process but rather an active engagement, a dance between our senses and our minds</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/><w:t>We perceive not merely with our senses but also with our minds, our memories, and our expectations</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Our past experiences, cultural background, and beliefs influence how we interpret sensory information</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> What we see, hear, smell, taste, and touch is filtered through the lens of our individual experiences, creating a unique and personal reality for each of us</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> This interplay between our senses, our minds, and our past experiences creates a dynamic and ever-changing perception of the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/><w:t>Perception is a complex process that defies simple definitions</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> It involves not only the physical reception of stimuli but also the cognitive interpretation of those stimuli</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> It is a process that is both subjective and objective, personal and universal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> It is the foundation of our understanding of the world and the key to unlocking the secrets of our consciousness</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/><w:t>Body:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/><w:t>Our perception is shaped by our senses</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> The eyes, ears, nose, tongue, and skin are the primary gateways through which we receive information about the world around us</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> These sensory organs convert physical stimuli into electrical signals that are transmitted to the brain for interpretation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> The brain then integrates these signals, creating a coherent and meaningful representation of the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> This process is incredibly complex and involves a multitude of neural pathways and brain regions working in concert</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:lastRenderedPageBreak/><w:br/><w:t>Perception is influenced by our past experiences</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Our memory plays a crucial role in shaping our perception of the present</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Past experiences create expectations about what we expect to see, hear, smell, taste, and touch</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> These expectations can influence how we interpret sensory information, leading to biased or inaccurate perceptions</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> For example, a person who has had a negative experience with a particular food may be more likely to perceive that food as being unappealing, even if it is not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:br/><w:t>Perception is also influenced by our beliefs and expectations</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Our cultural background, beliefs, and expectations can shape how we perceive the world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> For example, a person who believes in a particular religious doctrine may interpret ambiguous sensory information in a way that confirms their beliefs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> Similarly, a person who expects to see a certain outcome may be more likely to perceive evidence that supports that outcome, even if the evidence is weak or nonexistent</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/><w:sz w:val=\"24\"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$introParagraph = $d.Paragraphs(5).Range\n$introParagraph.InsertXML($introXml) | Out-Null\n\n# 5) Replace the Summary paragraph's text with its new \"Perception\" themed content.\n$summaryXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t>Perception is a complex and multifaceted process that involves the interaction of our senses, our minds, and our past experiences</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\"> It is a process that is both subjective and objective, personal and universal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\"> Perception is the foundation of our understanding of the world and the key to unlocking the secrets of our consciousness</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t xml:space=\"preserve\"> By understanding the nature of perception, we can gain a deeper appreciation for the intricate beauty of our subjective realities and the common threads that bind us together as a human species</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Aptos\" w:hAnsi=\"Aptos\"/><w:color w:val=\"000000\"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$summaryCount = $d.Paragraphs.Count\n$summaryParagraph = $d.Paragraphs($summaryCount).Range\n$summaryParagraph.InsertXML($summaryXml) | Out-Null\n\n# 6) A new, empty trailing paragraph is added at the very end of the document body.\n$lastCount = $d.Paragraphs.Count\n$d.Paragraphs($lastCount).Range.InsertParagraphAfter() | Out-Null\n"}
